# Auto-generated Excel COM-interop script implementing the 'Updated IPS AIP hipo turnover' commit.
# Updates Professional Voluntary Turnover CVD/YTD figures (and related monthly/quarterly rollups)
# across all seven IPS division/segment sheets, and clears the now-blank CVD figure on Integration.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("IPS Clutches & Brakes Division")
$ws.Range("D2").Value = 0.0425
$ws.Range("D3").Value = 0.0425
$ws.Range("D4").Value = 0.0425
$ws.Range("J4").Value = 0.0079
$ws.Range("L4").Value = 0.0079
$ws.Range("M4").Value = 0.0228
$ws.Range("N4").Value = 0.00708333333333333
$ws.Range("O4").Value = 0.00708333333333333
$ws.Range("P4").Value = 0.00708333333333333
$ws.Range("Q4").Value = 0.02125
$ws.Range("R4").Value = 0.00708333333333333
$ws.Range("S4").Value = 0.00708333333333333
$ws.Range("T4").Value = 0.00708333333333333
$ws.Range("U4").Value = 0.02125
$ws.Range("V4").Value = 0.085
$ws.Range("D5").Value = 0.48780487804878
$ws.Range("D6").Value = 0.48780487804878
$ws.Range("D7").Value = 0.48780487804878
$ws.Range("L7").Value = 0.75
$ws.Range("M7").Value = 0.5238
$ws.Range("N7").Value = 0.48780487804878
$ws.Range("O7").Value = 0.48780487804878
$ws.Range("P7").Value = 0.48780487804878
$ws.Range("Q7").Value = 0.48780487804878
$ws.Range("R7").Value = 0.48780487804878
$ws.Range("S7").Value = 0.48780487804878
$ws.Range("T7").Value = 0.48780487804878
$ws.Range("U7").Value = 0.48780487804878
$ws.Range("V7").Value = 0.48780487804878

$ws = $wb.Worksheets.Item("IPS Couplings Division")
$ws.Range("D2").Value = 0.0493
$ws.Range("D3").Value = 0.0493
$ws.Range("D4").Value = 0.0493
$ws.Range("G4").Value = 0.0035
$ws.Range("I4").Value = 0.0262
$ws.Range("K4").Value = 0.0085
$ws.Range("L4").Value = 0.0086
$ws.Range("M4").Value = 0.023
$ws.Range("N4").Value = 0.00821666666666667
$ws.Range("O4").Value = 0.00821666666666667
$ws.Range("P4").Value = 0.00821666666666667
$ws.Range("Q4").Value = 0.02465
$ws.Range("R4").Value = 0.00821666666666667
$ws.Range("S4").Value = 0.00821666666666667
$ws.Range("T4").Value = 0.00821666666666667
$ws.Range("U4").Value = 0.02465
$ws.Range("V4").Value = 0.0986
$ws.Range("D5").Value = 0.379310344827586
$ws.Range("D6").Value = 0.379310344827586
$ws.Range("D7").Value = 0.379310344827586
$ws.Range("H7").Value = 0.5
$ws.Range("I7").Value = 0.2667
$ws.Range("L7").Value = 0.625
$ws.Range("M7").Value = 0.5
$ws.Range("N7").Value = 0.379310344827586
$ws.Range("O7").Value = 0.379310344827586
$ws.Range("P7").Value = 0.379310344827586
$ws.Range("Q7").Value = 0.379310344827586
$ws.Range("R7").Value = 0.379310344827586
$ws.Range("S7").Value = 0.379310344827586
$ws.Range("T7").Value = 0.379310344827586
$ws.Range("U7").Value = 0.379310344827586
$ws.Range("V7").Value = 0.379310344827586

$ws = $wb.Worksheets.Item("IPS Gearing Division")
$ws.Range("D2").Value = 0.065
$ws.Range("D3").Value = 0.065
$ws.Range("D4").Value = 0.065
$ws.Range("J4").Value = 0.0114
$ws.Range("K4").Value = 0.0164
$ws.Range("L4").Value = 0.0165
$ws.Range("M4").Value = 0.0442
$ws.Range("N4").Value = 0.0108333333333333
$ws.Range("O4").Value = 0.0108333333333333
$ws.Range("P4").Value = 0.0108333333333333
$ws.Range("Q4").Value = 0.0325
$ws.Range("R4").Value = 0.0108333333333333
$ws.Range("S4").Value = 0.0108333333333333
$ws.Range("T4").Value = 0.0108333333333333
$ws.Range("U4").Value = 0.0325
$ws.Range("V4").Value = 0.13
$ws.Range("D5").Value = 0.48
$ws.Range("D6").Value = 0.48
$ws.Range("D7").Value = 0.48
$ws.Range("L7").Value = 0.4286
$ws.Range("M7").Value = 0.3846
$ws.Range("N7").Value = 0.48
$ws.Range("O7").Value = 0.48
$ws.Range("P7").Value = 0.48
$ws.Range("Q7").Value = 0.48
$ws.Range("R7").Value = 0.48
$ws.Range("S7").Value = 0.48
$ws.Range("T7").Value = 0.48
$ws.Range("U7").Value = 0.48
$ws.Range("V7").Value = 0.48

$ws = $wb.Worksheets.Item("IPS Industrial Components Divi")
$ws.Range("D2").Value = 0.0543
$ws.Range("D3").Value = 0.0543
$ws.Range("D4").Value = 0.0543
$ws.Range("F4").Value = 0.0108
$ws.Range("G4").Value = 0.0022
$ws.Range("I4").Value = 0.0216
$ws.Range("L4").Value = 0.0088
$ws.Range("M4").Value = 0.0327
$ws.Range("N4").Value = 0.00905
$ws.Range("O4").Value = 0.00905
$ws.Range("P4").Value = 0.00905
$ws.Range("Q4").Value = 0.02715
$ws.Range("R4").Value = 0.00905
$ws.Range("S4").Value = 0.00905
$ws.Range("T4").Value = 0.00905
$ws.Range("U4").Value = 0.02715
$ws.Range("V4").Value = 0.1086
$ws.Range("D5").Value = 0.52
$ws.Range("D6").Value = 0.52
$ws.Range("D7").Value = 0.52
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.4545
$ws.Range("N7").Value = 0.52
$ws.Range("O7").Value = 0.52
$ws.Range("P7").Value = 0.52
$ws.Range("Q7").Value = 0.52
$ws.Range("R7").Value = 0.52
$ws.Range("S7").Value = 0.52
$ws.Range("T7").Value = 0.52
$ws.Range("U7").Value = 0.52
$ws.Range("V7").Value = 0.52

$ws = $wb.Worksheets.Item("IPS Segment Functions")
$ws.Range("D2").Value = 0.048
$ws.Range("D3").Value = 0.048
$ws.Range("D4").Value = 0.048
$ws.Range("F4").Value = 0.0085
$ws.Range("G4").Value = 0.0042
$ws.Range("I4").Value = 0.0235
$ws.Range("J4").Value = 0.0068
$ws.Range("K4").Value = 0.0076
$ws.Range("L4").Value = 0.0101
$ws.Range("M4").Value = 0.0245
$ws.Range("N4").Value = 0.008
$ws.Range("O4").Value = 0.008
$ws.Range("P4").Value = 0.008
$ws.Range("Q4").Value = 0.024
$ws.Range("R4").Value = 0.008
$ws.Range("S4").Value = 0.008
$ws.Range("T4").Value = 0.008
$ws.Range("U4").Value = 0.024
$ws.Range("V4").Value = 0.096
$ws.Range("D5").Value = 0.638297872340426
$ws.Range("D6").Value = 0.638297872340426
$ws.Range("D7").Value = 0.638297872340426
$ws.Range("L7").Value = 0.6
$ws.Range("M7").Value = 0.5
$ws.Range("N7").Value = 0.638297872340426
$ws.Range("O7").Value = 0.638297872340426
$ws.Range("P7").Value = 0.638297872340426
$ws.Range("Q7").Value = 0.638297872340426
$ws.Range("R7").Value = 0.638297872340426
$ws.Range("S7").Value = 0.638297872340426
$ws.Range("T7").Value = 0.638297872340426
$ws.Range("U7").Value = 0.638297872340426
$ws.Range("V7").Value = 0.638297872340426

$ws = $wb.Worksheets.Item("Integration")
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()

$ws = $wb.Worksheets.Item("L1_IPS")
$ws.Range("D2").Value = 0.0502
$ws.Range("D3").Value = 0.0502
$ws.Range("D4").Value = 0.0502
$ws.Range("F4").Value = 0.0092
$ws.Range("I4").Value = 0.0225
$ws.Range("K4").Value = 0.0103
$ws.Range("L4").Value = 0.0101
$ws.Range("M4").Value = 0.0276
$ws.Range("N4").Value = 0.00836666666666667
$ws.Range("O4").Value = 0.00836666666666667
$ws.Range("P4").Value = 0.00836666666666667
$ws.Range("Q4").Value = 0.0251
$ws.Range("R4").Value = 0.00836666666666667
$ws.Range("S4").Value = 0.00836666666666667
$ws.Range("T4").Value = 0.00836666666666667
$ws.Range("U4").Value = 0.0251
$ws.Range("V4").Value = 0.1004
$ws.Range("D5").Value = 0.514970059880239
$ws.Range("D6").Value = 0.514970059880239
$ws.Range("D7").Value = 0.514970059880239
$ws.Range("H7").Value = 0.5294
$ws.Range("I7").Value = 0.5426
$ws.Range("L7").Value = 0.5667
$ws.Range("M7").Value = 0.4795
$ws.Range("N7").Value = 0.514970059880239
$ws.Range("O7").Value = 0.514970059880239
$ws.Range("P7").Value = 0.514970059880239
$ws.Range("Q7").Value = 0.514970059880239
$ws.Range("R7").Value = 0.514970059880239
$ws.Range("S7").Value = 0.514970059880239
$ws.Range("T7").Value = 0.514970059880239
$ws.Range("U7").Value = 0.514970059880239
$ws.Range("V7").Value = 0.514970059880239

Write-Output "Applied IPS AIP hipo turnover updates across all division/segment sheets."
